# Adds two new columns to the sheet:
#   I ("I0") - 1 for every data row except row 35, which is 5.
#   J ("IF") - mirrors the existing IP (column H) value for every data row
#              except row 35, which is 7 (H35 itself stays 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - same text style as the neighboring header cells (B1:H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$lastRow = 37

for ($r = 2; $r -le $lastRow; $r++) {
    $h = $ws.Cells.Item($r, 8).Value2

    if ($r -eq 35) {
        $i0 = 5
        $iF = 7
    } else {
        $i0 = 1
        $iF = $h
    }

    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $iF
}
